$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "CountTank" column to the table (expands table + autofilter A3:S21 -> A3:T21) ---
$tbl = $ws.ListObjects.Item(1)
$newCol = $tbl.ListColumns.Add()

# Header texts: row1 display header "攻城" then row3 technical column name "CountTank"
# (set in this order so shared-string indices come out 76="攻城", 77="CountTank")
$ws.Range("T1").Value = "攻城"
$ws.Range("T3").Value = "CountTank"

# Row2 type marker under the new column, same as the rest of row 2 ("int")
$ws.Range("T2").Value = "int"

# --- Column widths: columns D:T (4-20) all become width 3.75 ---
$wdCols = $ws.Range($ws.Cells.Item(1,4), $ws.Cells.Item(1,20))
$wdCols.ColumnWidth = 3.0

# --- New CountTank data values for rows 4-19 ---
$ws.Range("T4").Value = 0
$ws.Range("T5").Value = 4
$ws.Range("T6").Value = 0
$ws.Range("T7").Value = 0
$ws.Range("T8").Value = 0
$ws.Range("T9").Value = 0
$ws.Range("T10").Value = 0
$ws.Range("T11").Value = 0
$ws.Range("T12").Value = 0
$ws.Range("T13").Value = 0
$ws.Range("T14").Value = 0
$ws.Range("T15").Value = 0
$ws.Range("T16").Value = 0
$ws.Range("T17").Value = 0
$ws.Range("T18").Value = 0
$ws.Range("T19").Value = 0

# --- Fix other attr skill values on row 5 (机械) and row 9 (鸟) ---
$ws.Range("D5").Value = 0
$ws.Range("K5").Value = 4
$ws.Range("E9").Value = 3

# --- Highlight the corrected "skill" text cells with a yellow fill ---
$ws.Range("C5").Interior.Color = 65535
$ws.Range("C9").Interior.Color = 65535

# --- Extend the data-bar conditional formatting from S4:S19 to S4:T19 ---
$fc = $ws.Range("S4:S19").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("S4:T19"))

# --- Update the active selection to match the authored state ---
$ws.Range("K9").Select()
